$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every "Price" cell as literal text (it keeps
# European-style "." thousands separators and trailing zeros such as
# "1.00" or "26.30"). Plain Range.Value assignment auto-coerces any
# numeric-looking string to a real Number, which would both reformat
# the display (dropping trailing zeros) and change the stored cell type.
# So: force Text format on those specific cells first, assign the
# string, then reset the style back to "Normal" so the cell ends up
# with no extra style index (identical look to the untouched cells).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D11",
    "D12",
    "D15",
    "D16",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D25",
    "D26",
    "D28",
    "D29",
    "D31",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D48",
    "D49",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "43.951.45"
$ws.Range("E2").Value = "  +2.79%  "

# Row 3
$ws.Range("D3").Value = "2.268.94"
$ws.Range("E3").Value = "  +1.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "230.32"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("D6").Value = "0.635"
$ws.Range("E6").Value = "  +2.13%  "

# Row 7
$ws.Range("D7").Value = "63.35"
$ws.Range("E7").Value = "  +3.61%  "

# Row 8
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$ws.Range("D9").Value = "0.447"
$ws.Range("E9").Value = "  +9.81%  "

# Row 10
$ws.Range("E10").Value = "  +10.93%  "

# Row 11
$ws.Range("D11").Value = "56.94"
$ws.Range("E11").Value = "  -1.22%  "

# Row 12
$ws.Range("D12").Value = "26.30"
$ws.Range("E12").Value = "  +17.53%  "

# Row 13
$ws.Range("E13").Value = "  +2.18%  "

# Row 14
$ws.Range("D14").Value = "2.608.92"
$ws.Range("E14").Value = "  +1.53%  "

# Row 15
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  +0.31%  "

# Row 16
$ws.Range("D16").Value = "6.20"
$ws.Range("E16").Value = "  +9.18%  "

# Row 17
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +4.82%  "

# Row 18
$ws.Range("D18").Value = "2.265.17"
$ws.Range("E18").Value = "  +1.34%  "

# Row 19
$ws.Range("D19").Value = "43.901.89"
$ws.Range("E19").Value = "  +3.17%  "

# Row 20
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +6.92%  "

# Row 21
$ws.Range("D21").Value = "73.59"
$ws.Range("E21").Value = "  +1.44%  "

# Row 22
$ws.Range("D22").Value = "6.07"
$ws.Range("E22").Value = "  -2.08%  "

# Row 23
$ws.Range("D23").Value = "253.36"
$ws.Range("E23").Value = "  +3.16%  "

# Row 24
$ws.Range("E24").Value = "  +0.11%  "

# Row 25
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  -8.70%  "

# Row 26
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  -2.13%  "

# Row 28
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  +2.88%  "

# Row 29
$ws.Range("D29").Value = "171.70"
$ws.Range("E29").Value = "  +1.14%  "

# Row 30
$ws.Range("E30").Value = "  -3.36%  "

# Row 31
$ws.Range("D31").Value = "20.80"
$ws.Range("E31").Value = "  +1.69%  "

# Row 32
$ws.Range("E32").Value = "  -5.66%  "

# Row 33
$ws.Range("E33").Value = "  +3.05%  "

# Row 34
$ws.Range("D34").Value = "0.0697"
$ws.Range("E34").Value = "  +6.56%  "

# Row 35
$ws.Range("D35").Value = "4.79"
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("D36").Value = "4.90"
$ws.Range("E36").Value = "  -1.99%  "

# Row 37
$ws.Range("D37").Value = "3.81"
$ws.Range("E37").Value = "  +6.16%  "

# Row 38
$ws.Range("D38").Value = "6.53"
$ws.Range("E38").Value = "  +1.67%  "

# Row 39
$ws.Range("D39").Value = "2.31"
$ws.Range("E39").Value = "  -3.09%  "

# Row 40
$ws.Range("E40").Value = "  +3.46%  "

# Row 41
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "0.000236"
$ws.Range("E41").Value = "  +4.81%  "

# Row 42
$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "17.66"
$ws.Range("E43").Value = "  +8.10%  "

# Row 44
$ws.Range("E44").Value = "  +0.86%  "

# Row 45
$ws.Range("D45").Value = "8.25"
$ws.Range("E45").Value = "  -5.13%  "

# Row 46
$ws.Range("D46").Value = "97.94"
$ws.Range("E46").Value = "  +0.59%  "

# Row 47
$ws.Range("E47").Value = "  -1.21%  "

# Row 48
$ws.Range("D48").Value = "10.13"
$ws.Range("E48").Value = "  +10.99%  "

# Row 49
$ws.Range("D49").Value = "4.32"
$ws.Range("E49").Value = "  -1.91%  "

# Row 50
$ws.Range("D50").Value = "1.444.73"
$ws.Range("E50").Value = "  -1.03%  "

# Row 51
$ws.Range("D51").Value = "2.29"
$ws.Range("E51").Value = "  +2.71%  "

# Restore default styling on the cells we forced to Text format above
# so they end up indistinguishable (style-wise) from the untouched cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}